$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Y values in column E for rows 14-16
$ws.Range("E14").Value = 13
$ws.Range("E15").Value = 305
$ws.Range("E16").Value = 305

# Update the selected cell to E17
$ws.Range("E17").Select()
